$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 292 ("翼を広げなさい" post) entirely; all subsequent rows shift up by one.
$ws.Rows.Item(292).Delete()
